$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each literal value is written with a leading apostrophe so the engine
# keeps it typed as text (matching the source inline-string cells) instead
# of auto-coercing numeric-looking strings (e.g. "211.57") into numbers.
# ClearFormats() immediately afterwards drops the resulting quote-prefix
# style flag so the cell's style index is left untouched (still General/0).

$ws.Range("D2").Value = '''27.418.00'
$ws.Range("D2").ClearFormats()
$ws.Range("D3").Value = '''1.641.01'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '''  -1.55%  '
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = '''  +0.01%  '
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = '''211.57'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '''  -1.79%  '
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = '''0.533'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '''  +3.77%  '
$ws.Range("E6").ClearFormats()
$ws.Range("D8").Value = '''23.09'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '''  -2.18%  '
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = '''0.255'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '''  -3.18%  '
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = '''  -2.27%  '
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = '''0.0890'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '''  +0.93%  '
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = '''1.874.29'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '''  -1.44%  '
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = '''1.622.55'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '''  -2.52%  '
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = '''4.02'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '''  -3.09%  '
$ws.Range("E14").ClearFormats()
$ws.Range("E15").Value = '''  -0.56%  '
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = '''64.16'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '''  -2.97%  '
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = '''27.386.21'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '''  -0.83%  '
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = '''228.62'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '''  -6.53%  '
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = '''  -1.68%  '
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = '''7.47'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '''  -0.94%  '
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = '''0.999'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '''  -0.04%  '
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = '''4.31'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '''  -3.96%  '
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = '''9.28'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '''  -0.14%  '
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = '''  -1.16%  '
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = '''147.82'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '''  +1.04%  '
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = '''0.113'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '''  +1.47%  '
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = '''6.92'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '''  -3.77%  '
$ws.Range("E27").ClearFormats()
$ws.Range("E28").Value = '''  +0.04%  '
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = '''15.50'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '''  -5.35%  '
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = '''  -4.94%  '
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = '''0.0484'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '''  -4.06%  '
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = '''  -2.16%  '
$ws.Range("E32").ClearFormats()
$ws.Range("E33").Value = '''  -0.58%  '
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = '''1.411.39'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '''  -4.56%  '
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = '''1.56'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '''  -0.15%  '
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = '''  -0.12%  '
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = '''  -1.96%  '
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = '''0.878'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '''  -6.25%  '
$ws.Range("E38").ClearFormats()
$ws.Range("E40").Value = '''  +1.01%  '
$ws.Range("E40").ClearFormats()
$ws.Range("E42").Value = '''  -1.49%  '
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = '''5.45'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '''  +0.51%  '
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = '''2.22'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '''  +0.52%  '
$ws.Range("E44").ClearFormats()
$ws.Range("B45").Value = '''TrustWalletToken'
$ws.Range("B45").ClearFormats()
$ws.Range("C45").Value = '''https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("C45").ClearFormats()
$ws.Range("D45").Value = '''0.790'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '''  +0.25%  '
$ws.Range("E45").ClearFormats()
$ws.Range("B46").Value = '''Aave'
$ws.Range("B46").ClearFormats()
$ws.Range("C46").Value = '''https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("C46").ClearFormats()
$ws.Range("D46").Value = '''64.45'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '''  -7.48%  '
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = '''1.783.42'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '''  -1.42%  '
$ws.Range("E47").ClearFormats()
$ws.Range("E48").Value = '''  -4.14%  '
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = '''87.23'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '''  -2.35%  '
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = '''  -2.88%  '
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = '''0.0985'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '''  -3.82%  '
$ws.Range("E51").ClearFormats()
